$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure D and E columns are text-formatted so numeric-looking / percent-looking
# strings (prices, % changes) stay stored verbatim as text, matching the source data.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '256.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.27%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.58%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.644'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.27%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05919'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.60%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.603'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.54%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8563'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-1.45%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9115'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-3.38%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1376'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-1.84%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.04443'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '17.37%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06994'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.31%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03031'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-5.65%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09109'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-1.53%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001528'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-1.02%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0006032'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.30%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006105'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.57%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.474'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-1.15%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.133'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.84%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '0.20%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1288'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.50%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.895'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.26%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04224'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.23%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.35%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004615'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '8.25%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001200'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.07%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001716'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '13.74%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03798'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.31%'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006200'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-1.02%'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1098'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.25%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002312'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.51%'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '24.11%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005122'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-6.27%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.07%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05002'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-16.96%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '10,464.68%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.07%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.07%'
